# dang chinh lai doan luu vao excel
# Rework of the flashcard-data save routine: the old "StartTime" column (C)
# is dropped entirely and the remaining columns shift left. The column that
# used to hold "Step" (old D, date-only values) is regenerated with a fresh
# set of per-row dates and becomes the new column C.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "StartTime" column (C). Excel shifts D->C, E->D, F->E,
# G->F, H->G automatically, which also moves the "Step" header into C1 and
# carries each row's existing number formatting (date style) along with it.
$ws.Columns("C").Delete()

# Write the freshly computed "Step" dates (column C, rows 2-13).
$ws.Cells.Item(2, 3).Value = 44066
$ws.Cells.Item(3, 3).Value = 44067
$ws.Cells.Item(4, 3).Value = 44068
$ws.Cells.Item(5, 3).Value = 44069
$ws.Cells.Item(6, 3).Value = 44070
$ws.Cells.Item(7, 3).Value = 44071
$ws.Cells.Item(8, 3).Value = 44072
$ws.Cells.Item(9, 3).Value = 44073
$ws.Cells.Item(10, 3).Value = 44074
$ws.Cells.Item(11, 3).Value = 44075
$ws.Cells.Item(12, 3).Value = 44076
$ws.Cells.Item(13, 3).Value = 44066
